$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new Price value is a plain decimal (e.g. "247.73") must be
# explicitly formatted as Text first; otherwise Excel auto-converts the
# assigned string to a Number and can silently drop a trailing zero (as in
# "6.480" -> 6.48). The source cells are all text (inlineStr), so preserve that.

$ws.Range("D2").Value = "30.411.51"
$ws.Range("E2").Value = "  +0.25%  "

$ws.Range("D3").Value = "1.939.66"
$ws.Range("E3").Value = "  +0.20%  "

$ws.Range("E4").Value = "  +0.03%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7659"
$ws.Range("E5").Value = "  +7.88%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "247.73"
$ws.Range("E6").Value = "  -1.26%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "27.97"
$ws.Range("E8").Value = "  +1.47%  "

$ws.Range("E9").Value = "  -2.93%  "

$ws.Range("E10").Value = "  -2.62%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7842"
$ws.Range("E11").Value = "  -2.68%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08014"
$ws.Range("E12").Value = "  -0.70%  "

$ws.Range("D13").Value = "1.938.43"
$ws.Range("E13").Value = "  +0.14%  "

$ws.Range("E14").Value = "  -1.81%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "95.23"
$ws.Range("E15").Value = "  +0.65%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.56"
$ws.Range("E16").Value = "  -3.63%  "

$ws.Range("D17").Value = "30.411.08"
$ws.Range("E17").Value = "  +0.23%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "257.45"
$ws.Range("E18").Value = "  +1.70%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000008013"
$ws.Range("E19").Value = "  -2.34%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "5.856"
$ws.Range("E20").Value = "  +0.80%  "

$ws.Range("D21").Value = "2.195.47"
$ws.Range("E21").Value = "  +0.35%  "

$ws.Range("E22").Value = "  +0.02%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("E23").Value = "  -0.03%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.767"
$ws.Range("E24").Value = "  -3.41%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.615"
$ws.Range("E25").Value = "  -1.33%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "164.61"
$ws.Range("E26").Value = "  +0.31%  "

$ws.Range("E27").Value = "  -0.63%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1334"
$ws.Range("E28").Value = "  +3.63%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.297"
$ws.Range("E29").Value = "  -1.73%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.365"
$ws.Range("E30").Value = "  +1.18%  "

$ws.Range("E31").Value = "  -1.00%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.436"
$ws.Range("E32").Value = "  +0.48%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.152"
$ws.Range("E33").Value = "  -0.32%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05202"
$ws.Range("E34").Value = "  +0.41%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.281"
$ws.Range("E35").Value = "  +1.46%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7508"
$ws.Range("E36").Value = "  +0.53%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.779"
$ws.Range("E37").Value = "  +0.40%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01971"
$ws.Range("E38").Value = "  +0.01%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.808"
$ws.Range("E39").Value = "  +0.02%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "78.86"
$ws.Range("E40").Value = "  +0.03%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.480"
$ws.Range("E41").Value = "  +0.89%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.4525"
$ws.Range("E42").Value = "  -0.03%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.977"
$ws.Range("E43").Value = "  -1.84%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.002"
$ws.Range("E44").Value = "  +0.03%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8359"
$ws.Range("E45").Value = "  -1.31%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "101.42"
$ws.Range("E46").Value = "  -0.34%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.819"
$ws.Range("E47").Value = "  +0.89%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.543"
$ws.Range("E48").Value = "  +1.43%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "987.69"
$ws.Range("E49").Value = "  +11.31%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "37.37"
$ws.Range("E50").Value = "  +1.91%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.4168"
$ws.Range("E51").Value = "  -0.24%  "
